$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6481481481481481
$ws.Range("C2").Value = 0.7142857142857143
$ws.Range("D2").Value = 0.6796116504854369
$ws.Range("E2").Value = 49
$ws.Range("B3").Value = 0.6410256410256411
$ws.Range("C3").Value = 0.5681818181818182
$ws.Range("D3").Value = 0.6024096385542169
$ws.Range("E3").Value = 44
$ws.Range("B4").Value = 0.6451612903225806
$ws.Range("C4").Value = 0.6451612903225806
$ws.Range("D4").Value = 0.6451612903225806
$ws.Range("E4").Value = 0.6451612903225806
$ws.Range("B5").Value = 0.6445868945868947
$ws.Range("C5").Value = 0.6412337662337663
$ws.Range("D5").Value = 0.6410106445198269
$ws.Range("B6").Value = 0.6447783598321234
$ws.Range("C6").Value = 0.6451612903225806
$ws.Range("D6").Value = 0.6430859674212038
$ws.Range("B7").Value = 0.6428571428571429
$ws.Range("C7").Value = 0.7346938775510204
$ws.Range("D7").Value = 0.6857142857142857
$ws.Range("E7").Value = 49
$ws.Range("B8").Value = 0.6486486486486487
$ws.Range("C8").Value = 0.5454545454545454
$ws.Range("D8").Value = 0.5925925925925926
$ws.Range("E8").Value = 44
$ws.Range("B10").Value = 0.6457528957528957
$ws.Range("C10").Value = 0.6400742115027829
$ws.Range("D10").Value = 0.6391534391534391
$ws.Range("B11").Value = 0.6455972101133393
$ws.Range("D11").Value = 0.6416567104739147
$ws.Range("B12").Value = 0.660377358490566
$ws.Range("C12").Value = 0.7142857142857143
$ws.Range("D12").Value = 0.6862745098039216
$ws.Range("E12").Value = 49
$ws.Range("B13").Value = 0.65
$ws.Range("C13").Value = 0.5909090909090909
$ws.Range("D13").Value = 0.6190476190476191
$ws.Range("E13").Value = 44
$ws.Range("B14").Value = 0.6559139784946236
$ws.Range("C14").Value = 0.6559139784946236
$ws.Range("D14").Value = 0.6559139784946236
$ws.Range("E14").Value = 0.6559139784946236
$ws.Range("B15").Value = 0.655188679245283
$ws.Range("C15").Value = 0.6525974025974026
$ws.Range("D15").Value = 0.6526610644257703
$ws.Range("B16").Value = 0.6554676404950295
$ws.Range("C16").Value = 0.6559139784946236
$ws.Range("D16").Value = 0.6544682389084666
$ws.Range("B17").Value = 0.5686274509803921
$ws.Range("C17").Value = 0.5918367346938775
$ws.Range("D17").Value = 0.58
$ws.Range("E17").Value = 49
$ws.Range("B18").Value = 0.5238095238095238
$ws.Range("C18").Value = 0.5
$ws.Range("D18").Value = 0.5116279069767442
$ws.Range("E18").Value = 44
$ws.Range("B20").Value = 0.546218487394958
$ws.Range("C20").Value = 0.5459183673469388
$ws.Range("D20").Value = 0.5458139534883721
$ws.Range("B21").Value = 0.5474232703834222
$ws.Range("D21").Value = 0.5476519129782446
$ws.Range("B22").Value = 0.6181818181818182
$ws.Range("C22").Value = 0.6938775510204082
$ws.Range("D22").Value = 0.6538461538461539
$ws.Range("E22").Value = 49
$ws.Range("B23").Value = 0.6052631578947368
$ws.Range("C23").Value = 0.5227272727272727
$ws.Range("D23").Value = 0.5609756097560975
$ws.Range("E23").Value = 44
$ws.Range("B24").Value = 0.6129032258064516
$ws.Range("C24").Value = 0.6129032258064516
$ws.Range("D24").Value = 0.6129032258064516
$ws.Range("E24").Value = 0.6129032258064516
$ws.Range("B25").Value = 0.6117224880382774
$ws.Range("C25").Value = 0.6083024118738405
$ws.Range("D25").Value = 0.6074108818011257
$ws.Range("B26").Value = 0.6120697638524463
$ws.Range("C26").Value = 0.6129032258064516
$ws.Range("D26").Value = 0.6099074018035465